$d = $word.ActiveDocument

function Replace-ExactText {
    param($doc, [string]$oldText, [string]$newText)
    $rng = $doc.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { return $false }

    $start = $rng.Start
    $end = $rng.End

    # Compute common prefix/suffix so only the differing middle span is
    # rewritten -- this avoids touching (and thereby coalescing) any
    # neighbouring run that Find/Replace would otherwise also "dirty".
    $minLen = [Math]::Min($oldText.Length, $newText.Length)
    $prefixLen = 0
    while ($prefixLen -lt $minLen -and $oldText[$prefixLen] -eq $newText[$prefixLen]) {
        $prefixLen++
    }
    $oldRemaining = $oldText.Length - $prefixLen
    $newRemaining = $newText.Length - $prefixLen
    $maxSuffix = [Math]::Min($oldRemaining, $newRemaining)
    $suffixLen = 0
    while ($suffixLen -lt $maxSuffix -and $oldText[$oldText.Length - 1 - $suffixLen] -eq $newText[$newText.Length - 1 - $suffixLen]) {
        $suffixLen++
    }

    $midStart = $start + $prefixLen
    $midEnd = $end - $suffixLen
    $midNew = $newText.Substring($prefixLen, $newText.Length - $prefixLen - $suffixLen)

    $midRange = $doc.Range($midStart, $midEnd)
    $midRange.Text = $midNew
    return $true
}

$result0 = Replace-ExactText $d "Créditos-trabalho: 2" "Créditos-trabalho: 1"
Write-Output "Replace 0: $result0"

$result1 = Replace-ExactText $d "Carga horária: 75 h" "Carga horária: 45 h"
Write-Output "Replace 1: $result1"

$result2 = Replace-ExactText $d "Ativação: 01/01/2018" "Ativação: 01/01/2025"
Write-Output "Replace 2: $result2"

$result3 = Replace-ExactText $d "Conduzir os alunos no desenvolvimento de um projeto de conclusão de curso sobre tema específico relacionado à engenharia química." "Conduzir os alunos no desenvolvimento de um projeto de conclusão de curso sobre tema específico relacionado às atribuições da profissão."
Write-Output "Replace 3: $result3"

$result4 = Replace-ExactText $d "Lead students to develop a course completion project on a specific topic related to chemical engineering." "Lead students in the development of a course completion project on a specific topic related to the profession's duties."
Write-Output "Replace 4: $result4"

$result5 = Replace-ExactText $d "1) Metodologia Cientifica. 2) Projeto de Monografia. 3) Métodos de Pesquisa. 4) Normas de um Projeto de Pesquisa. 5) Pesquisa em Bases de Dados Bibliográficos. 6) Organização de Referências Bibliográficas." "1) Metodologia Cientifica. 2) Projetos de Monografia e de Relatório Técnico – diferenças e semelhanças. 3) Métodos de Pesquisa. 4) Normas aplicáveis a monografia e a relatórios técnicos. 5) Pesquisa em Bases de Dados Bibliográficos. 6) Organização de Referências Bibliográficas."
Write-Output "Replace 5: $result5"

$result6 = Replace-ExactText $d "1) Scientific Methodology. 2) Monography Project. 3) Research Methods. 4) Guidelines of a research project. 5) Research in Bibliographic Databases. 6) Organization of Bibliographic References." "1) Scientific Methodology. 2) Monograph and Technical Report Projects – differences and similarities. 3) Research Methods. 4) Standards applicable to monographs and technical reports. 5) Search in Bibliographic Databases. 6) Organization of Bibliographic References."
Write-Output "Replace 6: $result6"

$result7 = Replace-ExactText $d "1 – Metodologia Cientifica: Concepção e definição. 2 – Monografia Cientifica: O que é um projeto de pesquisa. As etapas de um projeto de pesquisa. 3 – Métodos de Pesquisa utilizados na Engenharia Química. 4- Pontos essenciais de um projeto de TCC. 5 – Pontos essenciais de uma monografia de TCC. 6 – Normas para elaboração de do texto e das Referências Bibliográficas. 7 – Mecanismos de busca em Bases de Dados tais como Web of Science, Scopus, Science Direct, etc. 8 – O uso do aplicativo Mendeley como forma de organizar e formatar as referências bibliográficas na monografia." "1. Metodologia Científica: Estudo das bases teóricas e práticas da pesquisa científica, abordando os métodos e técnicas para o desenvolvimento de investigações e análises científicas.2. Projetos de Monografia e de Relatório Técnico – Diferenças e Semelhanças: Análise comparativa entre a monografia acadêmica e o relatório técnico, destacando os aspectos comuns e os diferenciadores no contexto de projetos de engenharia.3. Métodos de Pesquisa: Exploração dos diferentes métodos de pesquisa aplicáveis em engenharia, incluindo pesquisa qualitativa, quantitativa, mista, experimental e bibliográfica.4. Normas Aplicáveis a Monografias e a Relatórios Técnicos: Discussão sobre as normas técnicas e acadêmicas vigentes que regulamentam a redação e formatação de monografias e relatórios técnicos.5. Pesquisa em Bases de Dados Bibliográficos: Instruções sobre como utilizar bases de dados bibliográficos para a coleta de informações e dados relevantes para o desenvolvimento de projetos de pesquisa em engenharia.6. Organização de Referências Bibliográficas: Orientações sobre a organização e citação de referências bibliográficas de acordo com as normas estabelecidas, utilizando ferramentas e softwares específicos para gestão bibliográfica."
Write-Output "Replace 7: $result7"

$result8 = Replace-ExactText $d "1 - Scientific Methodology: design and definition. 2 - Scientific Monography: What is a research project. The steps of a research project. 3 - Research Methods used in Chemical Engineering. 4 Key points of a Course Conclusion Paper project. 5 - Key points of a Course Conclusion Paper monography. 6 - Standards for preparation of the text and the references. 7 - Search engines in Databases such as Web of Science, Scopus, Science Direct, etc. 8 - The use of the Mendeley application as a way of organizing and formatting the bibliographic references in the monograph." "1. Scientific Methodology: Study of the theoretical and practical bases of scientific research, addressing methods and techniques for the development of scientific investigations and analyses.2. Monograph and Technical Report Projects – Differences and Similarities: Comparative analysis between the academic monograph and the technical report, highlighting the common and differentiating aspects in the context of engineering projects.3. Research Methods: Exploration of different research methods applicable in engineering, including qualitative, quantitative, mixed, experimental and bibliographic research.4. Standards Applicable to Monographs and Technical Reports: Discussion of current technical and academic standards that regulate the writing and formatting of monographs and technical reports.5. Search in Bibliographic Databases: Instructions on how to use bibliographic databases to collect information and data relevant to the development of engineering research projects.6. Organization of Bibliographic References: Guidance on the organization and citation of bibliographic references in accordance with established standards, using specific tools and software for bibliographic management."
Write-Output "Replace 8: $result8"

$result9 = Replace-ExactText $d "NASCIMENTO, L. P. Elaboração de Projetos de Pesquisa, Cengage Learning, 2012.SANTOS, C. R. Trabalho de Conclusão de Curso – Guia de elaboração passo a passo, Cengage Learning, 2010.BOOTH, W.; COLOMB, G.; WILLIAMS, J. A arte da Pesquisa. 3 ed. Martins Fontes. São Paulo. 2005. GIL, A.C. Como elaborar projetos de pesquisa. 5ed. Atlas, São Paulo, 2010." "BOOTH, W.; COLOMB, G.; WILLIAMS, J. A arte da Pesquisa. 3 ed. Martins Fontes. São Paulo. 2005. GIL, A.C. Como elaborar projetos de pesquisa. 5ed. Atlas, São Paulo, 2010.MIGUEL, PAULO CAUCHICK. Metodologia Científica Para Engenharia. Elsevier Brasil, 2014.NASCIMENTO, L. P. Elaboração de Projetos de Pesquisa, Cengage Learning, 2012.SANTOS, C. R. Trabalho de Conclusão de Curso – Guia de elaboração passo a passo, Cengage Learning, 2010."
Write-Output "Replace 9: $result9"
